$d = $word.ActiveDocument

# Change 1: "Content Type, Representation)" -> "Content Type, State / Representation)"
$d.Content.Find.Execute(
    "Content Type, Representation)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Content Type, State / Representation)", 2)

# Change 2: append clarification after "getOccurrences(Resource, Resource, Resource)?"
$d.Content.Find.Execute(
    "getOccurrences(Resource, Resource, Resource)?", $true, $false, $false, $false, $false,
    $true, 1, $false, "getOccurrences(Resource, Resource, Resource)? (CPPE / RCV / Kinds / Alignment schema / instances inference. Filtering / traversal).", 2)
